# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape, per commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 585
$ws1.Range("F4").Value = 44
$ws1.Range("F5").Value = 18
$ws1.Range("F7").Value = 14641
$ws1.Range("F9").Value = 667
$ws1.Range("F10").Value = 15156
$ws1.Range("F11").Value = 29
$ws1.Range("F12").Value = 8592
$ws1.Range("F13").Value = 313
$ws1.Range("F15").Value = 59
$ws1.Range("F16").Value = 175
$ws1.Range("F17").Value = 427
$ws1.Range("F18").Value = 181
$ws1.Range("F20").Value = 6
$ws1.Range("F24").Value = 1069
$ws1.Range("F27").Value = 52
$ws1.Range("F30").Value = 414
$ws1.Range("F31").Value = 20
$ws1.Range("F37").Value = 5304

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 585
$ws4.Range("F4").Value = 44
$ws4.Range("F5").Value = 18
$ws4.Range("F7").Value = 14641
$ws4.Range("F9").Value = 667
$ws4.Range("F10").Value = 15156
$ws4.Range("F11").Value = 29
$ws4.Range("F12").Value = 8592
$ws4.Range("F13").Value = 313
$ws4.Range("F16").Value = 59
$ws4.Range("F17").Value = 175
$ws4.Range("F18").Value = 427
$ws4.Range("F19").Value = 181
$ws4.Range("F21").Value = 6
$ws4.Range("F25").Value = 1069
$ws4.Range("F28").Value = 52
$ws4.Range("F33").Value = 414
$ws4.Range("F34").Value = 20
$ws4.Range("F40").Value = 5304
